$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 8 data rows (old rows 2-9), shifting everything else up.
$ws.Range("A2:C9").EntireRow.Delete() | Out-Null

# Append 7 new data rows of accelerometer readings after the existing data
# (which now ends at row 14, since rows 10-22 shifted up to rows 2-14).
$newData = @(
    @(-11.7663733065128, -6.984511554241188, 3.816117525100708),
    @(-1.357394456863403, -10.09910678863525, 3.819830894470215),
    @(1.284981921315195, -13.9884957075119, -13.69542229175569),
    @(1.272318005561828, -9.928469419479365, -8.66020488739013),
    @(2.418770149350169, -4.946553826332086, -8.305895447731023),
    @(-0.9577411413192878, -7.853628158569351, -5.728095054626454),
    @(-5.118649840354919, -10.2695299386978, -1.659017741680144)
)

$startRow = 15
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
